$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 28: fill in new timesheet entry for 2024-01-23 (serial 45314) ---
# Copy formatting (incl. date number format) from the row above, then set the date value.
$ws.Range("A27").Copy()
$ws.Range("A28").PasteSpecial(-4122)
$ws.Range("A28").Value = 45314

$ws.Range("B28").Formula = "=(1/60)*(8)"
$ws.Range("C28").Formula = "=(1/60)*(3)"
$ws.Range("D28").Formula = "=(1/60)*(15)"
$ws.Range("F28").Formula = "=SUM(B28:E28)"

# --- Row 29: new timesheet entry row for 2024-01-24 (serial 45315) ---
$ws.Range("A27").Copy()
$ws.Range("A29").PasteSpecial(-4122)
$ws.Range("A29").Value = 45315

$ws.Range("B29").Formula = "=(1/60)*(0)"
$ws.Range("C29").Formula = "=(1/60)*(0)"
$ws.Range("D29").Formula = "=(1/60)*(0)"
$ws.Range("E29").Formula = "=(1/60)*(0)"
$ws.Range("F29").Formula = "=SUM(B29:E29)"

# --- Update the saved selection to C29 ---
$ws.Range("C29").Select()
